$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing data rows 2..32 down to 3..33,
# preserving their original formatting).
$ws.Rows.Item(2).Insert(-4121, 1)

# Populate the new row 2 with the new weekly price entry.
$ws.Cells.Item(2, 1).Value = 11
$ws.Cells.Item(2, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(2, 3).Value = "Bíobío"
$ws.Cells.Item(2, 4).Value = 44496
$ws.Cells.Item(2, 5).Value = 8
$ws.Cells.Item(2, 6).Value = 100112013
$ws.Cells.Item(2, 7).Value = "Alcachofa"
$ws.Cells.Item(2, 8).Value = "Madrigal"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 350
$ws.Cells.Item(2, 11).Value = 7000
$ws.Cells.Item(2, 12).Value = 7500
$ws.Cells.Item(2, 13).Value = 7214
$ws.Cells.Item(2, 14).Value = "`$/caja 40 unidades"
$ws.Cells.Item(2, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(2, 16).Value = 180
$ws.Cells.Item(2, 17).Value = 40
$ws.Cells.Item(2, 18).Value = "Hortaliza"

# The freshly inserted row picked up the bold/centered header formatting from
# row 1 on Insert — strip it back to the plain (unstyled) look used by every
# other data row, then restore just the date number format on column D.
$ws.Range("A2:R2").ClearFormats()
$ws.Cells.Item(2, 4).NumberFormat = $ws.Cells.Item(3, 4).NumberFormat
